# "Add BS on 1D array codes"
# Adds three new Binary-Search-on-1D-array problems to the
# "BinarySearch - 1D Arrays" sheet (rows 6-8), each with its problem name
# (hyperlinked), platform, level, status and date - mirroring the format
# of the existing rows 2-5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BinarySearch - 1D Arrays")
$ws.Activate()

# --- Row 6: Floor and Ceil in Sorted Array ---------------------------------
$ws.Range("B6").Value = "Floor and Ceil in Sorted Array"
$ws.Hyperlinks.Add($ws.Range("B6"), "https://www.codingninjas.com/studio/problems/ceiling-in-a-sorted-array_1825401?utm_source=striver&utm_medium=website&utm_campaign=a_zcoursetuf") | Out-Null
$ws.Range("E6").Value = "CN"
$ws.Range("F6").Value = "Easy"
$ws.Range("G6").Value = "Pass"
$ws.Range("H6").Value = 45474

# --- Row 7: First and Last Occurrences in Array -----------------------------
$ws.Range("B7").Value = "First and Last Occurrences in Array"
$ws.Hyperlinks.Add($ws.Range("B7"), "https://leetcode.com/problems/find-first-and-last-position-of-element-in-sorted-array/") | Out-Null
$ws.Range("E7").Value = "LC"
$ws.Range("F7").Value = "Easy"
$ws.Range("G7").Value = "Pass"
$ws.Range("H7").Value = 45474

# --- Row 8: Count Occurrences in Sorted Array -------------------------------
$ws.Range("B8").Value = "Count Occurrences in Sorted Array"
$ws.Hyperlinks.Add($ws.Range("B8"), "https://www.codingninjas.com/studio/problems/count-occurrences_668372?utm_source=striver&utm_medium=website&utm_campaign=a_zcoursetuf") | Out-Null
$ws.Range("E8").Value = "CN"
$ws.Range("F8").Value = "Read"
$ws.Range("G8").Value = "Pass"
$ws.Range("H8").Value = 45474

# Re-apply the same cell formatting (borders/fills/date format/etc.) used by
# the matching existing rows, since entering values above can reset it.
$ws.Range("A3:I3").Copy()
$ws.Range("A6:I6").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("A2:I2").Copy()
$ws.Range("A7:I7").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("A5:I5").Copy()
$ws.Range("A8:I8").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Application.CutCopyMode = $false

# Keep the page set to portrait orientation.
$ws.PageSetup.Orientation = [Microsoft.Office.Interop.Excel.XlPageOrientation]::xlPortrait

# Leave the selection on the last edited cell, like the author did.
$ws.Range("H8").Select()
